# EPBDS-13605 Performance Improvement: The expression in the return cell is
# executed, despite the fact that corresponded the rule was not matched.
#
# Update the expected invocation-count test data on the "Testing" sheet:
# rows 36, 81, 122, 167 (columns D:H) change from 3,6,3,3,3 to 1,1,1,1,1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing")

$rows = @(36, 81, 122, 167)
foreach ($r in $rows) {
    $ws.Range("D$r`:H$r").Value = 1
}

$ws.Range("D36").Select() | Out-Null
